$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("quiz")

# Update marking value (per-correct-answer marks) from 3 to 5
$ws.Range("B11").Value = 5

# Update total correct marks from 60 to 100
$ws.Range("B12").Value = 100

# Update the "Corr/total" summary text from 54/84 to 100/140
$ws.Range("E12").Value = "100/140"
